$wb = $excel.ActiveWorkbook
$ws = $wb.ActiveSheet

# Add new "Save" header column in H1, matching the style of the other header cells (e.g. G1)
$ws.Range("H1").Value = "Save"
$ws.Range("G1").Copy()
$ws.Range("H1").PasteSpecial(-4122)  # xlPasteFormats
$excel.CutCopyMode = $false

# Fill in the new column's data values
$ws.Range("H2").Value = 1
$ws.Range("H3").Value = 0
